$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 3
    4  = 3
    5  = 3
    6  = 3
    7  = 3
    8  = 3
    9  = 3
    10 = 9
    11 = 6
    12 = 3
    13 = 6
    14 = 5
    15 = 3
    16 = 3
    17 = 3
    18 = 6
    19 = 3
    20 = 3
    21 = 6
    22 = 3
    23 = 3
    24 = 5
    25 = 3
    26 = 3
    27 = 3
    28 = 3
    29 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}
